# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (fund holdings detail) right before
#    the existing "2022-Q3" sheet.
# 2. Insert a new top row in the "总计" (totals) summary sheet for the
#    2022-Q4 quarter, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" summary sheet - insert a new row 2 for 2022-Q4
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Copy the number formatting/style of the index column (A) down from the
# row just below (which kept the original "s=2" style) so the new row
# matches the existing look.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.15

# Re-sequence the 0-based index column for the rows that shifted down.
for ($r = 3; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# Step 2: brand-new "2022-Q4" worksheet with the fund-holdings detail
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)              # "2022-Q3", new sheet goes before it
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# $q3's index shifted by one once the new sheet was inserted before it.
$q3 = $wb.Worksheets.Item(3)

# Borrow formatting (bold header + border, centered index column) from the
# "2022-Q3" sheet so the new sheet matches the workbook's existing style.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("A2").Copy()
$q4.Range("A2:A9").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (e.g. leading-zero fund
# codes, "92.80" style strings with significant trailing zeros) - force
# text storage so they round-trip exactly instead of being coerced to
# numbers.
$q4.Range("B2:B9").NumberFormat = "@"
$q4.Range("D2:G9").NumberFormat = "@"

$rows = @(
    @(0, "015784", "中信建投中证1000指数增强A", "6.76", "89.78", "0.64", "0.0433", 7),
    @(1, "004194", "招商中证1000指数增强A",     "2.57", "94.27", "1.16", "0.0298", 3),
    @(2, "004195", "招商中证1000指数增强C",     "2.14", "94.27", "1.16", "0.0248", 3),
    @(3, "015785", "中信建投中证1000指数增强C", "2.40", "89.78", "0.64", "0.0154", 7),
    @(4, "015148", "华安中证1000指数增强A",     "1.42", "92.79", "0.78", "0.0111", 5),
    @(5, "015149", "华安中证1000指数增强C",     "1.08", "92.79", "0.78", "0.0084", 5),
    @(6, "003646", "创金合信中证1000指数增强A", "0.52", "92.80", "1.34", "0.0070", 6),
    @(7, "003647", "创金合信中证1000指数增强C", "0.43", "92.80", "1.34", "0.0058", 6)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Restore the original "active sheet" (the last tab, "2021-Q2") so the
# selection highlight doesn't stay on the freshly inserted sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
